# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Most rows only get their Price (D) and Volume(1h) (E) columns refreshed; rows
# 47/48 additionally swap their Coin/Link (B/C) content (EnergySwap <-> Frax
# traded ranking positions). All cells in this sheet are stored as text, so
# numeric-looking price strings are written with a leading apostrophe to stop
# Excel from auto-converting them to numbers (which would silently drop
# significant trailing/leading zeros, e.g. "22.00" -> 22, "1.350" -> 1.35).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.285.47"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.684.93"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'218.16"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "'0.5244"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("D7").Value = "'1.008"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "'0.2709"
$ws.Range("E8").Value = "  +1.88%  "
$ws.Range("D9").Value = "'0.06417"
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").Value = "'22.00"
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").Value = "1.708.55"
$ws.Range("E12").Value = "  +2.28%  "
$ws.Range("D13").Value = "'4.561"
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "'0.5793"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "'0.000008451"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").Value = "'64.26"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "26.329.25"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").Value = "'4.923"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "'1.008"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("E20").Value = "  +0.29%  "
$ws.Range("D21").Value = "'188.52"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "'6.184"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").Value = "'1.009"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'144.47"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").Value = "'7.700"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("D26").Value = "'0.1233"
$ws.Range("E26").Value = "  +4.86%  "
$ws.Range("D27").Value = "'15.80"
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("D28").Value = "'0.06667"
$ws.Range("E28").Value = "  +13.62%  "
$ws.Range("D29").Value = "'1.350"
$ws.Range("E29").Value = "  +6.15%  "
$ws.Range("D30").Value = "'1.327"
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("D31").Value = "'3.576"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").Value = "'3.567"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("D33").Value = "'1.660"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("E34").Value = "  +1.47%  "
$ws.Range("D35").Value = "'0.6199"
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("D36").Value = "'2.398"
$ws.Range("E36").Value = "  +1.58%  "
$ws.Range("D37").Value = "'2.694"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("D38").Value = "'6.387"
$ws.Range("E38").Value = "  +5.84%  "
$ws.Range("D39").Value = "1.105.60"
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("D40").Value = "'0.01616"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").Value = "'0.8767"
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").Value = "'1.015"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "'100.72"
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("D44").Value = "1.832.23"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").Value = "'56.73"
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'8.191"
$ws.Range("E47").Value = "  +1.61%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.009"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("D49").Value = "'0.05272"
$ws.Range("E49").Value = "  +1.64%  "
$ws.Range("D50").Value = "'0.4306"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "'6.036"
$ws.Range("E51").Value = "  +2.90%  "
